$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.288.77'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '2.248.83'
$ws.Range('E3').Value = '  -1.37%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '247.08'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.54%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.629'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.19%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '74.53'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('E8').Value = '  +0.09%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.622'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -3.78%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '42.23'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +6.74%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0946'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -3.34%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '7.13'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -4.12%  '
$ws.Range('E13').Value = '  -3.35%  '
$ws.Range('D14').Value = '2.584.79'
$ws.Range('E14').Value = '  -1.28%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '14.55'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = '2.253.63'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').Value = '42.130.80'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = '0.0₃0981'
$ws.Range('E19').Value = '  -1.78%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '6.13'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.24%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '71.96'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E22').Value = '  +4.06%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '231.88'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -1.95%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '8.61'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +35.56%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.02%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '3.59'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -7.15%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -3.88%  '
$ws.Range('E29').Value = '  +1.01%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '169.90'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.39%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '20.72'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -1.58%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.0822'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -5.71%  '
$ws.Range('E33').Value = '  -5.38%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '30.33'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -4.71%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.124'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -2.50%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '4.57'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.59%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '5.07'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +6.42%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0306'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +0.31%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '13.45'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('E42').Value = '  -2.22%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '61.29'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.54%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '108.10'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +2.33%  '
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -4.14%  '
$ws.Range('E49').Value = '  -0.75%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.08%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '4.12'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -2.24%  '
